# Update "paises.xlsx" (countries / provincias Spain) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 12:03"

# --- Plain numeric refreshes (no reordering) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5306851
$ws.Range("C4").Value = 894
$ws.Range("D4").Value = 2756107
$ws.Range("E4").Value = 2382983
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = 167761

# Row 14: Iran
$ws.Range("B14").Value = 333699
$ws.Range("C14").Value = 2510
$ws.Range("D14").Value = 290244
$ws.Range("E14").Value = 24467
$ws.Range("G14").Value = 188
$ws.Range("H14").Value = 18988

# Row 18: Banglades
$ws.Range("B18").Value = 266498
$ws.Range("C18").Value = 2995
$ws.Range("D18").Value = 153089
$ws.Range("E18").Value = 109896
$ws.Range("G18").Value = 42
$ws.Range("H18").Value = 3513

# Row 33: Israel
$ws.Range("B33").Value = 87173
$ws.Range("C33").Value = 580
$ws.Range("D33").Value = 61577
$ws.Range("E33").Value = 24963
$ws.Range("G33").Value = 11
$ws.Range("H33").Value = 633

# Row 37: Oman
$ws.Range("B37").Value = 82299
$ws.Range("C37").Value = 249
$ws.Range("D37").Value = 77072
$ws.Range("E37").Value = 4688
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 539

# Row 73: El Salvador
$ws.Range("B73").Value = 21644
$ws.Range("C73").Value = 375
$ws.Range("D73").Value = 10056
$ws.Range("E73").Value = 11011

# Row 87: Consejo Danes para los Refugiados
$ws.Range("B87").Value = 9538
$ws.Range("C87").Value = 39
$ws.Range("D87").Value = 8421
$ws.Range("E87").Value = 892

# Row 88: Malasia
$ws.Range("B88").Value = 9114
$ws.Range("C88").Value = 11
$ws.Range("D88").Value = 8817
$ws.Range("E88").Value = 172

# Row 95: Finlandia
$ws.Range("B95").Value = 7642
$ws.Range("C95").Value = 19
$ws.Range("E95").Value = 329

# Row 111: Hong Kong
$ws.Range("B111").Value = 4244
$ws.Range("C111").Value = 62
$ws.Range("E111").Value = 1129
$ws.Range("G111").Value = 5
$ws.Range("H111").Value = 63

# Row 139: Nueva Zelanda
$ws.Range("B139").Value = 1579
$ws.Range("C139").Value = 9
$ws.Range("D139").Value = 1531
$ws.Range("E139").Value = 26

# --- Reorder: "Estado de Palestina" moves above "Dinamarca" (rows 77/78 swap),
#     and the Estado de Palestina figures get refreshed while Dinamarca's stay put ---
$ws.Range("A77").Value = "Estado de Palestina"
$ws.Range("B77").Value = 15184
$ws.Range("C77").Value = 309
$ws.Range("D77").Value = 8369
$ws.Range("E77").Value = 6711
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 104

$ws.Range("A78").Value = "Dinamarca"
$ws.Range("B78").Value = 14959
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 12988
$ws.Range("E78").Value = 1350
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 621

# --- Reorder: "Montserrat" moves above "Islas Malvinas" (rows 213/214 swap);
#     underlying figures for each country are unchanged ---
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
